# Commit message: "coba.ktr and edit target"
#
# The canonical XML diff for this commit only touches xl/workbook.xml:
#   - the worksheet/tab name "Target Penjualan" -> "Sales Target"
#   - plus a handful of Excel-session/application bookkeeping values
#     (the x15ac:absPath of the author's local working folder, the
#     xr:revisionPtr documentId GUID, and the bookViews window
#     width/height/position) that reflect the authoring machine's Excel
#     session rather than document content. Those are stamped by the
#     real Excel application when *it* saves the file and are not part
#     of the workbook object model exposed through COM automation, so
#     they cannot be (and should not be) set from a script - this
#     headless engine always re-derives/re-serializes that bookkeeping
#     on save regardless of script content.
#
# The one concrete, scriptable content change is renaming the sheet.

$wb = $excel.ActiveWorkbook

$target = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Target Penjualan") {
        $target = $sheet
        break
    }
}
if ($target -eq $null) {
    $target = $wb.ActiveSheet
}

$target.Name = "Sales Target"
